# Helper: write a full 16-column (A:P) row of data into a worksheet.
# Columns I (9), M (13), N (14) are numeric; everything else is forced to
# text (matching the source workbook's inlineStr convention) using the
# leading-apostrophe text-coercion trick. $null entries are skipped
# (left blank) except where an explicit empty string is required.
function Set-RowData($ws, $row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 1
        $val = $values[$i]
        if ($val -eq $null) { continue }
        if ($col -eq 9 -or $col -eq 13 -or $col -eq 14) {
            $ws.Cells.Item($row, $col).Value = $val
        } else {
            $ws.Cells.Item($row, $col).Value = "'" + $val
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# General sheet: two provider reassignments (AYKO -> NEW) plus three
# brand-new incident rows appended at the bottom (373-375).
# ---------------------------------------------------------------------
$general = $wb.Worksheets.Item("General")
$general.Cells.Item(225, 6).Value = "'NEW"
$general.Cells.Item(368, 6).Value = "'NEW"

Set-RowData $general 373 @("-541", "8/1/2025", "AYACUCHO 241", "3", "808663880", "PEBCOM", "Pendiente", "Colocar columna para pedir traspaso de nodo", 1, "Cambio", "Nodo Teco", "Pasante", -58.395015, -34.606755, "Almagro", "Capital Sur")
Set-RowData $general 374 @("-542", "8/1/2025", "Cramer 2141", "13", "808663881", "Sin Asignar", "Pendiente", "Cambiar columna 114 base corroida ", 1, "Cambio", "Sin equipos", "Pasante", -58.461582, -34.564296, "Saavedra", "Capital Norte")
Set-RowData $general 375 @("-543", "8/1/2025", "Pedro Ignacio Rivera 3258", "13", "", "NEW", "Pendiente", "Desmontar poste en desuso", 1, "Desmonte", "Sin equipos", "Poste", -58.46967, -34.561676, "Colegiales", "Capital Norte")

# ---------------------------------------------------------------------
# PEBCOM sheet: append the new "-541" record.
# ---------------------------------------------------------------------
$pebcom = $wb.Worksheets.Item("PEBCOM")
Set-RowData $pebcom 81 @("-541", "8/1/2025", "AYACUCHO 241", "3", "808663880", "PEBCOM", "Pendiente", "Colocar columna para pedir traspaso de nodo", 1, "Cambio", "Nodo Teco", "Pasante", -58.395015, -34.606755, "Almagro", "Capital Sur")

# ---------------------------------------------------------------------
# Sin_Asignar sheet: append the new "-542" record.
# ---------------------------------------------------------------------
$sinAsignar = $wb.Worksheets.Item("Sin_Asignar")
Set-RowData $sinAsignar 5 @("-542", "8/1/2025", "Cramer 2141", "13", "808663881", "Sin Asignar", "Pendiente", "Cambiar columna 114 base corroida ", 1, "Cambio", "Sin equipos", "Pasante", -58.461582, -34.564296, "Saavedra", "Capital Norte")

# ---------------------------------------------------------------------
# AYKO sheet: two records ("-406" Olof palme 4144, "-536" Olof palme
# 4142) are reassigned away to the NEW sheet, so remove them here.
# Row 88 shifts to row 87 once row 62 is removed.
# ---------------------------------------------------------------------
$ayko = $wb.Worksheets.Item("AYKO")
$ayko.Rows.Item(62).Delete()
$ayko.Rows.Item(87).Delete()

# ---------------------------------------------------------------------
# NEW sheet: receives the two reassigned AYKO records (now tagged as
# "NEW") inserted in their chronological slots, plus the brand-new
# "-543" record appended at the end.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Item("NEW")
$newSheet.Rows.Item(35).Insert()
$newSheet.Rows.Item(65).Insert()

Set-RowData $newSheet 35 @("-406", "5/8/2025", "Olof palme 4144", "12", "805791925", "NEW", "Pendiente", "Tensar 2 riendas a pique columna 168", 1, "Tensor", "Sin equipos", "Terminal", -58.488252, -34.553391, "Saavedra", "Capital Norte")
Set-RowData $newSheet 65 @("-536", "7/29/2025", "Olof palme 4142", "12", "ICD30249764 ", "NEW", "Pendiente", "Aplomar o desmontar poste", 1, "Desmonte", "Sin equipos", "Poste", -58.488239, -34.55341, "Saavedra", "Capital Norte")
Set-RowData $newSheet 67 @("-543", "8/1/2025", "Pedro Ignacio Rivera 3258", "13", "", "NEW", "Pendiente", "Desmontar poste en desuso", 1, "Desmonte", "Sin equipos", "Poste", -58.46967, -34.561676, "Colegiales", "Capital Norte")
